# The deck ships with two DrawingML theme parts: theme1.xml (the "Office
# Theme" palette) and theme2.xml (the "Integral" palette - the one actually
# wired to the single slide master / every slide, via
# slideMaster1.xml.rels). This edit swaps the two themes' color schemes so
# the slides' theme becomes the "Office" palette (matching what used to be
# theme1.xml's clrScheme); the font/format schemes are already identical
# between the two theme parts and are left untouched.
#
# PowerPoint's ThemeColorScheme collection (reached off a Slide / the active
# design) is keyed 1-12 in the standard DrawingML clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
# and .RGB uses the usual COM RGB(r,g,b) = r + g*256 + b*65536 packing.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
